# Apln-Aplnr NATMI LR-pairs sheet: update with new TPM-derived statistics.
# Each changed cell is written with its new value from the updated
# (TPM-normalized) computation; unchanged cells are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 36.70383333333333
$ws.Range("H2").Value = 110.1115
$ws.Range("I2").Value = 0.7793963013403331
$ws.Range("J2").Value = 0.78322622712988
$ws.Range("M2").Value = 97.999779
$ws.Range("N2").Value = 293.999337
$ws.Range("O2").Value = 0.954667399411267
$ws.Range("P2").Value = 0.9562924695444477
$ws.Range("Q2").Value = 3596.9675551195
$ws.Range("R2").Value = 32372.7079960755
$ws.Range("S2").Value = 0.744064240111336
$ws.Range("T2").Value = 0.7489933429540134

# Row 3
$ws.Range("G3").Value = 36.70383333333333
$ws.Range("H3").Value = 110.1115
$ws.Range("I3").Value = 0.7793963013403331
$ws.Range("J3").Value = 0.78322622712988
$ws.Range("O3").Value = 0.005525698806503622
$ws.Range("P3").Value = 0.00553510485524996
$ws.Range("Q3").Value = 20.81956432011111
$ws.Range("R3").Value = 187.376078881
$ws.Range("S3").Value = 0.004306709212109616
$ws.Range("T3").Value = 0.004335239292545707

# Row 4
$ws.Range("G4").Value = 36.70383333333333
$ws.Range("H4").Value = 110.1115
$ws.Range("I4").Value = 0.7793963013403331
$ws.Range("J4").Value = 0.78322622712988
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.238873
$ws.Range("N4").Value = 9.716619
$ws.Range("O4").Value = 0.03155156568193249
$ws.Range("P4").Value = 0.03160527392322827
$ws.Range("Q4").Value = 118.8790547798333
$ws.Range("R4").Value = 1069.9114930185
$ws.Range("S4").Value = 0.02459117359399476
$ws.Range("T4").Value = 0.02475407945229646

# Row 5
$ws.Range("G5").Value = 36.70383333333333
$ws.Range("H5").Value = 110.1115
$ws.Range("I5").Value = 0.7793963013403331
$ws.Range("J5").Value = 0.78322622712988
$ws.Range("M5").Value = 0.5233300000000001
$ws.Range("N5").Value = 1.04666
$ws.Range("O5").Value = 0.005098032824481149
$ws.Range("P5").Value = 0.003404473922923818
$ws.Range("Q5").Value = 19.20821709833333
$ws.Range("R5").Value = 115.24930259
$ws.Range("S5").Value = 0.003973387927512219
$ws.Range("T5").Value = 0.002666473266013684

# Row 6
$ws.Range("G6").Value = 36.70383333333333
$ws.Range("H6").Value = 110.1115
$ws.Range("I6").Value = 0.7793963013403331
$ws.Range("J6").Value = 0.78322622712988
$ws.Range("M6").Value = 0.3241076666666667
$ws.Range("N6").Value = 0.972323
$ws.Range("O6").Value = 0.003157303275815759
$ws.Range("P6").Value = 0.003162677754150398
$ws.Range("Q6").Value = 11.89599377938889
$ws.Range("R6").Value = 107.0639440145
$ws.Range("S6").Value = 0.00246079049538052
$ws.Range("T6").Value = 0.002477092165010819

# Row 7
$ws.Range("I7").Value = 0.2059338919414416
$ws.Range("J7").Value = 0.2069458437845949
$ws.Range("M7").Value = 97.999779
$ws.Range("N7").Value = 293.999337
$ws.Range("O7").Value = 0.954667399411267
$ws.Range("P7").Value = 0.9562924695444477
$ws.Range("Q7").Value = 950.3990800815951
$ws.Range("R7").Value = 8553.591720734357
$ws.Range("S7").Value = 0.1965983730703769
$ws.Range("T7").Value = 0.1979007520147298

# Row 8
$ws.Range("I8").Value = 0.2059338919414416
$ws.Range("J8").Value = 0.2069458437845949
$ws.Range("O8").Value = 0.005525698806503622
$ws.Range("P8").Value = 0.00553510485524996
$ws.Range("S8").Value = 0.00113792866091947
$ws.Range("T8").Value = 0.001145466944705911

# Row 9
$ws.Range("I9").Value = 0.2059338919414416
$ws.Range("J9").Value = 0.2069458437845949
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.238873
$ws.Range("N9").Value = 9.716619
$ws.Range("O9").Value = 0.03155156568193249
$ws.Range("P9").Value = 0.03160527392322827
$ws.Range("Q9").Value = 31.41049858593167
$ws.Range("R9").Value = 282.694487273385
$ws.Range("S9").Value = 0.006497536717726381
$ws.Range("T9").Value = 0.00654058008008573

# Row 10
$ws.Range("I10").Value = 0.2059338919414416
$ws.Range("J10").Value = 0.2069458437845949
$ws.Range("M10").Value = 0.5233300000000001
$ws.Range("N10").Value = 1.04666
$ws.Range("O10").Value = 0.005098032824481149
$ws.Range("P10").Value = 0.003404473922923818
$ws.Range("Q10").Value = 5.075239512316668
$ws.Range("R10").Value = 30.45143707390001
$ws.Range("S10").Value = 0.001049857740790623
$ws.Range("T10").Value = 0.0007045417286221196

# Row 11
$ws.Range("I11").Value = 0.2059338919414416
$ws.Range("J11").Value = 0.2069458437845949
$ws.Range("M11").Value = 0.3241076666666667
$ws.Range("N11").Value = 0.972323
$ws.Range("O11").Value = 0.003157303275815759
$ws.Range("P11").Value = 0.003162677754150398
$ws.Range("Q11").Value = 3.143186968282778
$ws.Range("R11").Value = 28.288682714545
$ws.Range("S11").Value = 0.0006501957516282021
$ws.Range("T11").Value = 0.0006545030164514218

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.69084
$ws.Range("H12").Value = 1.38168
$ws.Range("I12").Value = 0.01466980671822532
$ws.Range("J12").Value = 0.009827929085525242
$ws.Range("M12").Value = 97.999779
$ws.Range("N12").Value = 293.999337
$ws.Range("O12").Value = 0.954667399411267
$ws.Range("P12").Value = 0.9562924695444477
$ws.Range("Q12").Value = 67.70216732436
$ws.Range("R12").Value = 406.21300394616
$ws.Range("S12").Value = 0.0140047862295541
$ws.Range("T12").Value = 0.00939837457570464

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.69084
$ws.Range("H13").Value = 1.38168
$ws.Range("I13").Value = 0.01466980671822532
$ws.Range("J13").Value = 0.009827929085525242
$ws.Range("O13").Value = 0.005525698806503622
$ws.Range("P13").Value = 0.00553510485524996
$ws.Range("Q13").Value = 0.39186609432
$ws.Range("R13").Value = 2.35119656592
$ws.Range("S13").Value = 0.00008106093347453648
$ws.Range("T13").Value = 0.00005439861799834307

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.69084
$ws.Range("H14").Value = 1.38168
$ws.Range("I14").Value = 0.01466980671822532
$ws.Range("J14").Value = 0.009827929085525242
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 3.238873
$ws.Range("N14").Value = 9.716619
$ws.Range("O14").Value = 0.03155156568193249
$ws.Range("P14").Value = 0.03160527392322827
$ws.Range("Q14").Value = 2.23754302332
$ws.Range("R14").Value = 13.42525813992
$ws.Range("S14").Value = 0.0004628553702113406
$ws.Range("T14").Value = 0.0003106143908460877

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.69084
$ws.Range("H15").Value = 1.38168
$ws.Range("I15").Value = 0.01466980671822532
$ws.Range("J15").Value = 0.009827929085525242
$ws.Range("M15").Value = 0.5233300000000001
$ws.Range("N15").Value = 1.04666
$ws.Range("O15").Value = 0.005098032824481149
$ws.Range("P15").Value = 0.003404473922923818
$ws.Range("Q15").Value = 0.3615372972000001
$ws.Range("R15").Value = 1.4461491888
$ws.Range("S15").Value = 0.00007478715617830677
$ws.Range("T15").Value = 0.00003345892828801522

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.69084
$ws.Range("H16").Value = 1.38168
$ws.Range("I16").Value = 0.01466980671822532
$ws.Range("J16").Value = 0.009827929085525242
$ws.Range("M16").Value = 0.3241076666666667
$ws.Range("N16").Value = 0.972323
$ws.Range("O16").Value = 0.003157303275815759
$ws.Range("P16").Value = 0.003162677754150398
$ws.Range("Q16").Value = 0.22390654044
$ws.Range("R16").Value = 1.34343924264
$ws.Range("S16").Value = 0.00004631702880703684
$ws.Range("T16").Value = 0.00003108257268815835
